$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# Insert a new column before column A, shifting the existing headers
# (F9, F9PL, ...) one column to the right.
$ws.Columns("A:A").Insert()

# New header for the inserted column. Match the formatting already
# applied to the rest of row 1 (bold, centered, bordered).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Value = "DIA"

# New data rows under the header row.
$row2 = @(7, $null, $null, 4, $null, $null, $null, $null, $null, $null, 5, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 3, 5, $null, 5, $null, $null, $null, $null, $null)
$row3 = @(8, 5, $null, 4, $null, $null, 5, $null, $null, $null, 3, $null, $null, $null, $null, 5, $null, $null, $null, 3, $null, $null, $null, 5, 5, $null, $null, $null, $null, $null, $null, $null)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $cell = $ws.Cells.Item(2, $i + 1)
    if ($null -ne $row2[$i]) {
        $cell.Value = $row2[$i]
    }
}

for ($i = 0; $i -lt $row3.Length; $i++) {
    $cell = $ws.Cells.Item(3, $i + 1)
    if ($null -ne $row3[$i]) {
        $cell.Value = $row3[$i]
    }
}
